$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.153.07"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "1.579.00"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D5").Value = "209.44"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "0.497"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.801.27"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "1.575.44"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "64.41"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "26.161.32"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "7.25"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "207.67"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "144.44"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "6.98"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "15.20"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "1.276.11"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "0.613"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "5.56"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "0.764"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").Value = "62.37"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "1.714.72"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "88.92"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("E51").Value = "  -0.09%  "
